$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -5
    3  = 1
    4  = 5
    5  = -2
    6  = 2
    7  = -8
    8  = 3
    9  = -3
    10 = 3
    11 = 4
    12 = -2
    13 = 1
    14 = 2
    15 = -3
    17 = 6
    18 = -2
    19 = -1
    20 = -1
    21 = -6
    22 = -1
    23 = 1
    24 = 2
    25 = 1
    26 = -5
    27 = -1
    28 = -2
    29 = -2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
